$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "2024-12-25 19:53:07"
$ws.Range("F2").Value = 1000
$ws.Range("G2").Value = 200
$ws.Range("H2").Value = 173
